$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing service code WS-PS-01 -> WS-PS-01-2 (row 24, ServiceCode column D) ---
$ws.Range("D24").Value = "WS-PS-01-2"

# --- Add new row 36: "Image Upload" service (used for image upload in Add Vehicle screen) ---

# Bring over the row's cell formatting (borders/fills/number-formats) from row 35,
# which has the same overall column layout/styles we need for the new row.
$ws.Range("B35:L35").Copy()
$ws.Range("B36:L36").PasteSpecial(-4122)

$ws.Range("B36").Value = "Image Upload"
$ws.Range("C36").Value = "Image Upload"
$ws.Range("D36").Value = "WS-IMG-01"

# Column F ("logActivity") stores the literal text "false" (not a boolean) in this sheet.
# Copy it from another row that already holds that literal text/format so it lands
# as text re-using the existing shared string instead of being auto-typed to a boolean.
$ws.Range("F24").Copy()
$ws.Range("F36").PasteSpecial(-4163)

$ws.Range("G36").Value = "upload"
$ws.Range("H36").Value = "/image"
$ws.Range("I36").Value = "POST"

$ws.Range("M36").Formula = '=_xlfn.CONCAT("INSERT INTO ",CHAR(34),"M_CTL_CONFIG",CHAR(34)," VALUES(''",D36,"'',''CONNON_CONFIG'', 0, ''",C36,"'', ''{}'', 0, 0, CURRENT_TIMESTAMP, ''ATUL'', null, null);")'
$ws.Range("N36").Formula = '=_xlfn.CONCAT(IF(I36="GET","@GetMapping(",IF(I36="POST","@PostMapping(",IF(I36="DELETE","@DeleteMapping(",IF(I36="PUT","@PutMapping(","")))),CHAR(34),H36,CHAR(34),")")'
$ws.Range("O36").Formula = '=_xlfn.CONCAT("@ServiceInfo(serviceCode = ",CHAR(34),D36,,CHAR(34),", serviceName = ",CHAR(34),C36,CHAR(34), ", queryId = ",CHAR(34),E36,CHAR(34),", logActivity =",F36,")")'

$wb.Application.Calculate()

# --- Reflect the post-edit UI scroll/selection state ---
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("G41").Select()
